# Updates the cryptocurrency price (column D) and 1h volume change (column E)
# figures on Sheet1, mirroring the latest scrape performed by the GitHub
# Actions job that refreshes cryptos.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = "<new price text>"; E = "<new volume text>" }
# (D is only present for rows whose price actually changed)
$updates = @{
  2  = @{ D = "56.961.40";  E = "  +1.30%  " }
  3  = @{ D = "3.244.84";   E = "  +0.36%  " }
  4  = @{              E = "  -0.02%  " }
  5  = @{ D = "397.02";     E = "  -0.47%  " }
  6  = @{ D = "108.12";     E = "  -2.63%  " }
  7  = @{              E = "  +5.00%  " }
  8  = @{              E = "  +0.01%  " }
  9  = @{ D = "0.620";      E = "  -1.05%  " }
  10 = @{ D = "39.33";      E = "  -0.88%  " }
  11 = @{              E = "  +6.88%  " }
  12 = @{              E = "  +1.95%  " }
  13 = @{ D = "3.753.16";   E = "  +0.23%  " }
  14 = @{ D = "8.32";       E = "  +2.77%  " }
  15 = @{              E = "  -1.26%  " }
  16 = @{ D = "3.230.05";   E = "  -0.43%  " }
  17 = @{              E = "  -3.50%  " }
  18 = @{ D = "11.07";      E = "  +4.91%  " }
  19 = @{ D = "56.811.53";  E = "  +1.04%  " }
  20 = @{              E = "  -0.82%  " }
  21 = @{              E = "  +9.96%  " }
  22 = @{ D = "13.03";      E = "  -0.82%  " }
  23 = @{ D = "294.04";     E = "  +1.62%  " }
  24 = @{ D = "74.43";      E = "  +0.12%  " }
  25 = @{              E = "  -1.67%  " }
  26 = @{ D = "28.13";      E = "  -0.59%  " }
  27 = @{              E = "  -0.64%  " }
  28 = @{ D = "7.81";       E = "  -4.76%  " }
  29 = @{              E = "  -1.38%  " }
  30 = @{ D = "7.24";       E = "  -4.19%  " }
  31 = @{              E = "  +0.02%  " }
  32 = @{ D = "42.07";      E = "  +13.64%  " }
  33 = @{ D = "11.17";      E = "  -1.09%  " }
  34 = @{              E = "  -2.52%  " }
  35 = @{ D = "0.0485";     E = "  -2.65%  " }
  36 = @{              E = "  +1.49%  " }
  37 = @{ D = "51.38";      E = "  +0.40%  " }
  38 = @{              E = "  -0.10%  " }
  39 = @{ D = "3.46";       E = "  -3.43%  " }
  40 = @{              E = "  -3.43%  " }
  41 = @{ D = "136.71";     E = "  -1.62%  " }
  42 = @{              E = "  +3.07%  " }
  43 = @{              E = "  -1.97%  " }
  44 = @{              E = "  -2.64%  " }
  45 = @{ D = "16.82";      E = "  +0.24%  " }
  46 = @{              E = "  -3.57%  " }
  47 = @{              E = "  +8.61%  " }
  48 = @{ D = "22.51";      E = "  -0.24%  " }
  49 = @{ D = "2.153.65";   E = "  +0.94%  " }
  50 = @{              E = "  -5.83%  " }
  51 = @{ D = "1.97";       E = "  -6.90%  " }
}

foreach ($row in $updates.Keys) {
  $entry = $updates[$row]

  if ($entry.ContainsKey("D")) {
    $price = $entry.D
    # Plain decimal-looking strings (e.g. "397.02") would otherwise be
    # auto-coerced into numbers (dropping significant trailing zeros), so
    # force them to stay text with a leading quote prefix. Values that use
    # dots as thousands separators (e.g. "56.961.40") are never parsed as
    # numbers by Excel, so they can be assigned as-is.
    if ($price -match '^[0-9]+\.[0-9]+$') {
      $ws.Cells.Item($row, 4).Value = "'" + $price
    } else {
      $ws.Cells.Item($row, 4).Value = $price
    }
  }

  if ($entry.ContainsKey("E")) {
    $ws.Cells.Item($row, 5).Value = $entry.E
  }
}
